# Update Sheets via scheduled runner: apply latest market-price snapshot
# to the Odin_Profits leve-profit calculations across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 3754.7778
$ws.Range("I5").Value = 1080.7
$ws.Range("K5").Value = 1080.7
$ws.Range("M5").Value = -965.7

$ws.Range("H38").Value = 633.9091
$ws.Range("J38").Value = 1925
$ws.Range("L38").Value = 5775
$ws.Range("N38").Value = -6519

$ws.Range("H62").Value = 6763.478
$ws.Range("I62").Value = 6293.1577
$ws.Range("K62").Value = 6293.1577
$ws.Range("M62").Value = -5669.1577

$ws.Range("H65").Value = 6763.478
$ws.Range("I65").Value = 6293.1577
$ws.Range("K65").Value = 31465.7885
$ws.Range("M65").Value = -28345.7885

$ws.Range("H86").Value = 111112870
$ws.Range("J86").Value = 2029
$ws.Range("L86").Value = 2029
$ws.Range("N86").Value = -4275

$ws.Range("H89").Value = 111112870
$ws.Range("J89").Value = 2029
$ws.Range("L89").Value = 10145
$ws.Range("N89").Value = -21377

$ws.Range("H106").Value = 1695.6666
$ws.Range("I106").Value = 1790.1818
$ws.Range("J106").Value = 1279.8
$ws.Range("K106").Value = 1790.1818
$ws.Range("L106").Value = 1279.8
$ws.Range("M106").Value = -1159.1818
$ws.Range("N106").Value = -2541.8

$ws.Range("H132").Value = 8082.25
$ws.Range("I132").Value = 4182.8423
$ws.Range("K132").Value = 12548.5269
$ws.Range("M132").Value = -10018.5269

$ws.Range("H138").Value = 3964.1714
$ws.Range("I138").Value = 1003.1667
$ws.Range("K138").Value = 3009.5001
$ws.Range("M138").Value = 2130.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5157.077
$ws.Range("I2").Value = 4828.15
$ws.Range("K2").Value = 4828.15
$ws.Range("M2").Value = -4715.15

$ws.Range("H13").Value = 14287485
$ws.Range("J13").Value = 2478.6
$ws.Range("L13").Value = 2478.6
$ws.Range("N13").Value = -2766.6

$ws.Range("H14").Value = 16667330
$ws.Range("I14").Value = 16667330
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 16667330
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -16667155
$ws.Range("N14").ClearContents()

$ws.Range("H45").Value = 2015.05
$ws.Range("I45").Value = 1914.875
$ws.Range("K45").Value = 1914.875
$ws.Range("M45").Value = -1537.875

$ws.Range("H74").Value = 2234.8
$ws.Range("I74").Value = 2474.6667
$ws.Range("K74").Value = 2474.6667
$ws.Range("M74").Value = -1600.6667

$ws.Range("H77").Value = 2234.8
$ws.Range("I77").Value = 2474.6667
$ws.Range("K77").Value = 12373.3335
$ws.Range("M77").Value = -8005.333500000001

$ws.Range("H116").Value = 5157.077
$ws.Range("I116").Value = 4828.15
$ws.Range("K116").Value = 4828.15
$ws.Range("M116").Value = -2534.15

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5157.077
$ws.Range("I3").Value = 4828.15
$ws.Range("K3").Value = 4828.15
$ws.Range("M3").Value = -4714.15

$ws.Range("H22").Value = 727.6667
$ws.Range("I22").Value = 446.4
$ws.Range("K22").Value = 446.4
$ws.Range("M22").Value = -273.4

$ws.Range("H64").Value = 5563.857
$ws.Range("I64").Value = 1766.6666
$ws.Range("J64").Value = 8411.75
$ws.Range("K64").Value = 1766.6666
$ws.Range("L64").Value = 8411.75
$ws.Range("M64").Value = -1541.6666
$ws.Range("N64").Value = -8861.75

$ws.Range("H67").Value = 5563.857
$ws.Range("I67").Value = 1766.6666
$ws.Range("J67").Value = 8411.75
$ws.Range("K67").Value = 1766.6666
$ws.Range("L67").Value = 8411.75
$ws.Range("M67").Value = -986.6666
$ws.Range("N67").Value = -9971.75

$ws.Range("H80").Value = 14535.667
$ws.Range("J80").Value = 19247.154
$ws.Range("L80").Value = 19247.154
$ws.Range("N80").Value = -21243.154

$ws.Range("H83").Value = 14535.667
$ws.Range("J83").Value = 19247.154
$ws.Range("L83").Value = 96235.76999999999
$ws.Range("N83").Value = -106219.77

$ws.Range("H105").Value = 2277.8147
$ws.Range("I105").Value = 2315.9546
$ws.Range("K105").Value = 2315.9546
$ws.Range("M105").Value = -568.9546

$ws.Range("H134").Value = 8155.857
$ws.Range("I134").Value = 4233.7144
$ws.Range("K134").Value = 12701.1432
$ws.Range("M134").Value = -10166.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2143
$ws.Range("I6").Value = 2143
$ws.Range("K6").Value = 2143
$ws.Range("M6").Value = -2030

$ws.Range("H22").Value = 1865.2354
$ws.Range("I22").Value = 599.2222
$ws.Range("K22").Value = 599.2222
$ws.Range("M22").Value = -249.2222

$ws.Range("H99").Value = 5773.727
$ws.Range("J99").Value = 7557.857
$ws.Range("L99").Value = 7557.857
$ws.Range("N99").Value = -10553.857

$ws.Range("H126").Value = 5773.727
$ws.Range("J126").Value = 7557.857
$ws.Range("L126").Value = 22673.571
$ws.Range("N126").Value = -27613.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 23151234
$ws.Range("I140").Value = 28738256
$ws.Range("J140").Value = 4999.4287
$ws.Range("K140").Value = 86214768
$ws.Range("L140").Value = 14998.2861
$ws.Range("M140").Value = -86209588
$ws.Range("N140").Value = -25358.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 9826.182000000001
$ws.Range("I113").Value = 1859
$ws.Range("K113").Value = 1859
$ws.Range("M113").Value = 311

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6447.9165
$ws.Range("I122").Value = 4750
$ws.Range("K122").Value = 14250
$ws.Range("M122").Value = -11800

$ws.Range("H132").Value = 2683.2856
$ws.Range("I132").Value = 2372.4443
$ws.Range("K132").Value = 7117.3329
$ws.Range("M132").Value = -4587.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 613.1429000000001
$ws.Range("I107").Value = 465.33334
$ws.Range("K107").Value = 1396.00002
$ws.Range("M107").Value = 523.9999800000001

$ws.Range("H122").Value = 5547
$ws.Range("J122").Value = 5500
$ws.Range("L122").Value = 16500
$ws.Range("N122").Value = -21400

$ws.Range("H132").Value = 3476.4167
$ws.Range("I132").Value = 2098.7407
$ws.Range("K132").Value = 6296.222099999999
$ws.Range("M132").Value = -3766.222099999999

$ws.Range("H136").Value = 59819.11
$ws.Range("I136").Value = 73735
$ws.Range("J136").Value = 11113.5
$ws.Range("K136").Value = 221205
$ws.Range("L136").Value = 33340.5
$ws.Range("M136").Value = -218655
$ws.Range("N136").Value = -38440.5
